# update scripts wuth new tpm
# Recomputed NATMI LR-pair statistics (Efnb2-Epha4) for sheet1 with updated TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 45.71598933333333
$ws.Range("H2").Value = 137.147968
$ws.Range("I2").Value = 0.6549002937372808
$ws.Range("J2").Value = 0.6549002937372808
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.23061133333333
$ws.Range("N2").Value = 30.691834
$ws.Range("O2").Value = 0.4855635428718841
$ws.Range("P2").Value = 0.4855635428718841
$ws.Range("Q2").Value = 467.7025185881458
$ws.Range("R2").Value = 4209.322667293311
$ws.Range("S2").Value = 0.3179957068549116
$ws.Range("T2").Value = 0.3179957068549116
$ws.Range("G3").Value = 45.71598933333333
$ws.Range("H3").Value = 137.147968
$ws.Range("I3").Value = 0.6549002937372808
$ws.Range("J3").Value = 0.6549002937372808
$ws.Range("O3").Value = 0.4164864079521221
$ws.Range("P3").Value = 0.4164864079521222
$ws.Range("Q3").Value = 401.1663248126791
$ws.Range("R3").Value = 3610.496923314111
$ws.Range("S3").Value = 0.2727570709054297
$ws.Range("T3").Value = 0.2727570709054297
$ws.Range("G4").Value = 45.71598933333333
$ws.Range("H4").Value = 137.147968
$ws.Range("I4").Value = 0.6549002937372808
$ws.Range("J4").Value = 0.6549002937372808
$ws.Range("M4").Value = 2.034752
$ws.Range("N4").Value = 6.104255999999999
$ws.Range("O4").Value = 0.09657305490303886
$ws.Range("P4").Value = 0.09657305490303887
$ws.Range("Q4").Value = 93.02070072797865
$ws.Range("R4").Value = 837.1863065518079
$ws.Range("S4").Value = 0.06324572202310669
$ws.Range("T4").Value = 0.0632457220231067
$ws.Range("G5").Value = 45.71598933333333
$ws.Range("H5").Value = 137.147968
$ws.Range("I5").Value = 0.6549002937372808
$ws.Range("J5").Value = 0.6549002937372808
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.02901266666666667
$ws.Range("N5").Value = 0.087038
$ws.Range("O5").Value = 0.001376994272954919
$ws.Range("P5").Value = 0.001376994272954919
$ws.Range("Q5").Value = 1.326342759864889
$ws.Range("R5").Value = 11.937084838784
$ws.Range("S5").Value = 0.0009017939538327295
$ws.Range("T5").Value = 0.0009017939538327295
$ws.Range("I6").Value = 0.1818108415648851
$ws.Range("J6").Value = 0.1818108415648851
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.23061133333333
$ws.Range("N6").Value = 30.691834
$ws.Range("O6").Value = 0.4855635428718841
$ws.Range("P6").Value = 0.4855635428718841
$ws.Range("Q6").Value = 129.8417321227207
$ws.Range("R6").Value = 1168.575589104486
$ws.Range("S6").Value = 0.0882807163627644
$ws.Range("T6").Value = 0.08828071636276438
$ws.Range("I7").Value = 0.1818108415648851
$ws.Range("J7").Value = 0.1818108415648851
$ws.Range("O7").Value = 0.4164864079521221
$ws.Range("P7").Value = 0.4164864079521222
$ws.Range("S7").Value = 0.07572174433011136
$ws.Range("T7").Value = 0.07572174433011136
$ws.Range("I8").Value = 0.1818108415648851
$ws.Range("J8").Value = 0.1818108415648851
$ws.Range("M8").Value = 2.034752
$ws.Range("N8").Value = 6.104255999999999
$ws.Range("O8").Value = 0.09657305490303886
$ws.Range("P8").Value = 0.09657305490303887
$ws.Range("Q8").Value = 25.82404076473599
$ws.Range("R8").Value = 232.4163668826239
$ws.Range("S8").Value = 0.01755802838441334
$ws.Range("T8").Value = 0.01755802838441334
$ws.Range("I9").Value = 0.1818108415648851
$ws.Range("J9").Value = 0.1818108415648851
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.02901266666666667
$ws.Range("N9").Value = 0.087038
$ws.Range("O9").Value = 0.001376994272954919
$ws.Range("P9").Value = 0.001376994272954919
$ws.Range("Q9").Value = 0.3682140559113333
$ws.Range("R9").Value = 3.313926503202
$ws.Range("S9").Value = 0.0002503524875959608
$ws.Range("T9").Value = 0.0002503524875959608
$ws.Range("G10").Value = 11.24784666666667
$ws.Range("H10").Value = 33.74354
$ws.Range("I10").Value = 0.161130015850732
$ws.Range("J10").Value = 0.161130015850732
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.23061133333333
$ws.Range("N10").Value = 30.691834
$ws.Range("O10").Value = 0.4855635428718841
$ws.Range("P10").Value = 0.4855635428718841
$ws.Range("Q10").Value = 115.0723475835956
$ws.Range("R10").Value = 1035.65112825236
$ws.Range("S10").Value = 0.07823886135948425
$ws.Range("T10").Value = 0.07823886135948425
$ws.Range("G11").Value = 11.24784666666667
$ws.Range("H11").Value = 33.74354
$ws.Range("I11").Value = 0.161130015850732
$ws.Range("J11").Value = 0.161130015850732
$ws.Range("O11").Value = 0.4164864079521221
$ws.Range("P11").Value = 0.4164864079521222
$ws.Range("Q11").Value = 98.7019503487622
$ws.Range("R11").Value = 888.3175531388598
$ws.Range("S11").Value = 0.06710846151493986
$ws.Range("T11").Value = 0.06710846151493986
$ws.Range("G12").Value = 11.24784666666667
$ws.Range("H12").Value = 33.74354
$ws.Range("I12").Value = 0.161130015850732
$ws.Range("J12").Value = 0.161130015850732
$ws.Range("M12").Value = 2.034752
$ws.Range("N12").Value = 6.104255999999999
$ws.Range("O12").Value = 0.09657305490303886
$ws.Range("P12").Value = 0.09657305490303887
$ws.Range("Q12").Value = 22.88657850069333
$ws.Range("R12").Value = 205.97920650624
$ws.Range("S12").Value = 0.01556081786728026
$ws.Range("T12").Value = 0.01556081786728026
$ws.Range("G13").Value = 11.24784666666667
$ws.Range("H13").Value = 33.74354
$ws.Range("I13").Value = 0.161130015850732
$ws.Range("J13").Value = 0.161130015850732
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.02901266666666667
$ws.Range("N13").Value = 0.087038
$ws.Range("O13").Value = 0.001376994272954919
$ws.Range("P13").Value = 0.001376994272954919
$ws.Range("Q13").Value = 0.3263300260577778
$ws.Range("R13").Value = 2.93697023452
$ws.Range("S13").Value = 0.0002218751090275932
$ws.Range("T13").Value = 0.0002218751090275932
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1507006666666667
$ws.Range("H14").Value = 0.452102
$ws.Range("I14").Value = 0.00215884884710222
$ws.Range("J14").Value = 0.00215884884710222
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 10.23061133333333
$ws.Range("N14").Value = 30.691834
$ws.Range("O14").Value = 0.4855635428718841
$ws.Range("P14").Value = 0.4855635428718841
$ws.Range("Q14").Value = 1.541759948340889
$ws.Range("R14").Value = 13.875839535068
$ws.Range("S14").Value = 0.001048258294723836
$ws.Range("T14").Value = 0.001048258294723836
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1507006666666667
$ws.Range("H15").Value = 0.452102
$ws.Range("I15").Value = 0.00215884884710222
$ws.Range("J15").Value = 0.00215884884710222
$ws.Range("O15").Value = 0.4164864079521221
$ws.Range("P15").Value = 0.4164864079521222
$ws.Range("Q15").Value = 1.322426430557555
$ws.Range("R15").Value = 11.901837875018
$ws.Range("S15").Value = 0.0008991312016411837
$ws.Range("T15").Value = 0.0008991312016411836
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1507006666666667
$ws.Range("H16").Value = 0.452102
$ws.Range("I16").Value = 0.00215884884710222
$ws.Range("J16").Value = 0.00215884884710222
$ws.Range("M16").Value = 2.034752
$ws.Range("N16").Value = 6.104255999999999
$ws.Range("O16").Value = 0.09657305490303886
$ws.Range("P16").Value = 0.09657305490303887
$ws.Range("Q16").Value = 0.3066384829013333
$ws.Range("R16").Value = 2.759746346112
$ws.Range("S16").Value = 0.0002084866282385648
$ws.Range("T16").Value = 0.0002084866282385648
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1507006666666667
$ws.Range("H17").Value = 0.452102
$ws.Range("I17").Value = 0.00215884884710222
$ws.Range("J17").Value = 0.00215884884710222
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.02901266666666667
$ws.Range("N17").Value = 0.087038
$ws.Range("O17").Value = 0.001376994272954919
$ws.Range("P17").Value = 0.001376994272954919
$ws.Range("Q17").Value = 0.004372228208444445
$ws.Range("R17").Value = 0.039350053876
$ws.Range("S17").Value = 0.000002972722498635085
$ws.Range("T17").Value = 0.000002972722498635085
